$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying weekly price records (rows 2-8) were reordered (re-sorted
# by date), so each destination row now carries the D/M/N/O/P/S values that
# used to belong to a different row. Apply the permutation explicitly.

$map = @{
    2 = 7
    3 = 4
    4 = 5
    5 = 3
    6 = 8
    7 = 2
    8 = 6
}

# Snapshot the original values before overwriting anything.
$orig = @{}
foreach ($r in 2..8) {
    $orig[$r] = @{
        D = $ws.Range("D$r").Value()
        M = $ws.Range("M$r").Value()
        N = $ws.Range("N$r").Value()
        O = $ws.Range("O$r").Value()
        P = $ws.Range("P$r").Value()
        S = $ws.Range("S$r").Value()
    }
}

foreach ($r in 2..8) {
    $src = $map[$r]
    $ws.Range("D$r").Value = $orig[$src].D
    $ws.Range("M$r").Value = $orig[$src].M
    $ws.Range("N$r").Value = $orig[$src].N
    $ws.Range("O$r").Value = $orig[$src].O
    $ws.Range("P$r").Value = $orig[$src].P
    $ws.Range("S$r").Value = $orig[$src].S
}
